$d = $word.ActiveDocument

# --- 1) Two empty paragraphs right before the date: spacing after 60 -> 0 twips (3pt -> 0pt) ---
$d.Paragraphs(4).SpaceAfter = 0
$d.Paragraphs(5).SpaceAfter = 0

# --- 2) Date paragraph text: "Dec 25 2018" -> "Jan 4 2019" (collapses to a single run) ---
$pDate = $d.Paragraphs(6)
$rDate = $pDate.Range
$rDate.End = $rDate.End - 1
$rDate.Text = "Jan 4 2019"

# --- 3) "Respawn Entertainment" address paragraph gains spacing after = 360 twips (18pt) ---
$d.Paragraphs(7).SpaceAfter = 18

# --- 4) "Dear Respawn," paragraph spacing after 240 -> 120 twips (12pt -> 6pt) ---
$d.Paragraphs(8).SpaceAfter = 6

# --- 5) "I am..." paragraph: "Gameplay Programmer" -> "Gameplay Software Engineer";
#        insert " (Titanfall)" before " position at Respawn!", with a _GoBack bookmark
#        wrapping the text between "Titanfall" and ")".
$p9 = $d.Paragraphs(9)
$r9 = $p9.Range
$r9.End = $r9.End - 1
$r9.Text = "I am Shantanu Mane, a Gameplay Software Engineer highly skilled in C++ and 3D Math. I have accrued experience with Player Input, Animations and AI working on Action Games. I am currently a graduate student at the University of Utah studying in the final semester for my EAE - Game Engineering Master" + [char]0x2019 + "s degree. And I am beyond thrilled to be applying for the Gameplay Software Engineer (Titanfall) position at Respawn!"

$p9b = $d.Paragraphs(9)
$r9b = $p9b.Range.Duplicate
$r9b.Find.Execute("Titanfall") | Out-Null
$bmPos = $r9b.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done-part-1"

# --- 6) Swap paragraphs 11 & 12 (with text edits) in one combined replace ---
$p11 = $d.Paragraphs(11)
$p12 = $d.Paragraphs(12)
$start = $p11.Range.Start
$end = $p12.Range.End
$end = $end - 1
$combined = $d.Range($start, $end)
$newPara11 = 'I love Action Games that bring out raw excitement and power and have mastery in skill! I happen to be someone with a lot of enthusiasm and energy. I’ve played Titanfall 2 and I can say it evoked these feelings in me gracefully and masterfully. I like parkour and I loved wall-running in the game. It being so fluid and natural to pull-off opened up avenues for me to explore all my abilities and get creative with them. I’ve had adrenaline-filled combat encounters in the game where I was on point, pushing my abilities to their extent, wall-running and shooting enemies then launching off and landing to go sliding into an enemy to finish them off with a powerful and bone-shattering melee strike!'
$newPara12 = 'Combat, weapons, action and animation are where my passion truly lies. I have taken and continue to take time to learn combat design and combat systems. I am also putting time into learning about animation programming and am working towards creating a gameplay animation system related to but not limited to combat. I feel the free-flowing movement of Titanfall fulfills the power fantasy of being a nimble and agile swashbuckler as the Pilot and then the powerful arsenal of weapons you have when using your Titan makes you go gung-ho, keep the trigger pulled and lay all your firepower into your enemies. I think the game does an excellent job of making difficult feats achievable with deftness. This puts the player in the power fantasy the game wants them to feel from the very first moment they wall-run and leaves room for creativity past that point, with the systems being intuitive and rewarding the player for engaging in them.'
$combined.Text = $newPara11 + [char]13 + $newPara12

# --- 7) Paragraph 13 text edit ---
$p13 = $d.Paragraphs(13)
$r13 = $p13.Range
$r13.End = $r13.End - 1
$r13.Text = 'Games with fluid, fast-paced action are what Respawn is the flag-bearer of! I like the studio’s focus on pushing the envelope for games in new directions and the emphasis placed on gameplay feel. I admire how you have created such innovative gameplay that is fluid and responsive and achieved a very cool power fantasy with it too! It is the kind of player experience that I would like to be creating myself. And I would absolutely love to be a part of Respawn Entertainment and play my part in creating games that are trailblazers of these creative fantasies!'

Write-Output "done-part-2"
